# Apply the STAMP Session-Reflector text revisions described in the diff.
$p = $ppt.ActivePresentation

# --- Slide 3: "Requirements, Goals and Scope" -------------------------------
# "Avoid provisioning and maintaining each test session on Session-Reflector - Stateless mode"
#   -> "Avoid provisioning and maintaining each test session on Session-Reflector"
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange
$tr3.Paragraphs(4, 1).Text = "Avoid provisioning and maintaining each test session on Session-Reflector"

# --- Slide 7: "STAMP Return Path Control Code Sub-TLV - Usage" --------------
# First bullet's "Stateless mode of STAMP Session-Reflector as defined in RFC 8762 ..."
#   -> "Avoid maintaining each test session (session id, source-address) on Session-Reflector"
# and the (now duplicate) bullet further down with that same text is removed.
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(2)
$tr7 = $sh7.TextFrame.TextRange
$tr7.Paragraphs(1, 1).Text = "Avoid maintaining each test session (session id, source-address) on Session-Reflector"
$tr7.Paragraphs(6, 1).Delete()

# --- Slide 8: "STAMP Return Path Control Code Sub-TLV" ----------------------
# "Stateless mode of STAMP Session-Reflector as defined in RFC 8762 - per each test session"
#   -> "Avoid maintaining each test session (session id, source-address) on Session-Reflector"
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(2)
$tr8 = $sh8.TextFrame.TextRange
$tr8.Paragraphs(1, 1).Text = "Avoid maintaining each test session (session id, source-address) on Session-Reflector"

# --- Slide 9: "STAMP Return Path Segment List Sub-TLVs - Usage" -------------
# Remove the leading bullet "Stateless mode of STAMP Session-Reflector as defined in RFC 8762 ..." entirely.
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(2)
$tr9 = $sh9.TextFrame.TextRange
$tr9.Paragraphs(1, 1).Delete()
